# Apply the edits captured in the commit:
#  1. Slide 16's table switches to a different built-in table style
#     ({268E6C19-95ED-4354-9776-FCE8A10C793B} -> {A0D7EA1D-1999-4976-8EA1-A0EAC0C34FC3})
#  2. The presentation's theme colour scheme (the "Integral" scheme used by
#     the slide master / theme1.xml) is swapped for the "Office Theme" scheme
#     that was otherwise only referenced by the notes master (theme2.xml).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$s = $p.Slides.Item(16)
$tableShape = $s.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{A0D7EA1D-1999-4976-8EA1-A0EAC0C34FC3}")

# --- 2. Theme colour scheme -------------------------------------------------
# Re-point the twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) from the "Integral" palette to the standard "Office"
# palette. RGB values must be supplied as BGR-packed OLE colour integers.
$officeTheme = @{
    1  = 0         # dk1      000000
    2  = 16777215  # lt1      FFFFFF
    3  = 6968388   # dk2      44546A
    4  = 15132391  # lt2      E7E6E6
    5  = 13998939  # accent1  5B9BD5
    6  = 3243501   # accent2  ED7D31
    7  = 10855845  # accent3  A5A5A5
    8  = 49407     # accent4  FFC000
    9  = 12874308  # accent5  4472C4
    10 = 4697456    # accent6  70AD47
    11 = 12673797   # hlink    0563C1
    12 = 7491477    # folHlink 954F72
}

$themeColors = $s.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeTheme[$i]
}
